$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.150009477396964
$ws.Range("D2").Value = 0.1057530028625955
$ws.Range("E2").Value = 0.176351281074016
$ws.Range("F2").Value = 2.667917174710539
$ws.Range("G2").Value = 1.924414947864392
$ws.Range("H2").Value = 1.644247430075907
$ws.Range("J2").Value = 0.2823616107578868
$ws.Range("K2").Value = 0.4989959832325894
$ws.Range("L2").Value = 0.1831972504997879
$ws.Range("B3").Value = 1.129709051907696
$ws.Range("D3").Value = 0.1058104823714352
$ws.Range("E3").Value = 0.1769547775959062
$ws.Range("F3").Value = 2.663717253031209
$ws.Range("G3").Value = 1.916178103679371
$ws.Range("H3").Value = 1.645702696212823
$ws.Range("J3").Value = 0.2835398909239224
$ws.Range("K3").Value = 0.4543013397459674
$ws.Range("L3").Value = 0.1755946966196831
$ws.Range("B4").Value = 1.117799055574636
$ws.Range("D4").Value = 0.1058689386586096
$ws.Range("E4").Value = 0.1773572025392971
$ws.Range("F4").Value = 2.662258506802431
$ws.Range("G4").Value = 1.912030327492957
$ws.Range("H4").Value = 1.647214342052735
$ws.Range("J4").Value = 0.2843094493501059
$ws.Range("K4").Value = 0.4270592133098319
$ws.Range("L4").Value = 0.1710233903785934
$ws.Range("B5").Value = 1.113085524910815
$ws.Range("D5").Value = 0.1058986099945649
$ws.Range("E5").Value = 0.1775292264870894
$ws.Range("F5").Value = 2.661945683103724
$ws.Range("G5").Value = 1.910568615098128
$ws.Range("H5").Value = 1.647985809803743
$ws.Range("J5").Value = 0.2846346570401668
$ws.Range("K5").Value = 0.4160085736265842
$ws.Range("L5").Value = 0.1691849537376413
$ws.Range("B6").Value = 1.112311308868698
$ws.Range("D6").Value = 0.1059038909643242
$ws.Range("E6").Value = 0.1775582765809061
$ws.Range("F6").Value = 2.661910748278885
$ws.Range("G6").Value = 1.910339696128602
$ws.Range("H6").Value = 1.648123301792822
$ws.Range("J6").Value = 0.2846893590187705
$ws.Range("K6").Value = 0.4141766977911914
$ws.Range("L6").Value = 0.168881159690784
$ws.Range("B7").Value = 1.117734920253014
$ws.Range("D7").Value = 0.1058693150963137
$ws.Range("E7").Value = 0.177359489967932
$ws.Range("F7").Value = 2.662253147677831
$ws.Range("G7").Value = 1.912009689257886
$ws.Range("H7").Value = 1.647224116850353
$ws.Range("J7").Value = 0.284313788198336
$ws.Range("K7").Value = 0.4269099746195479
$ws.Range("L7").Value = 0.1709984976274228
$ws.Range("B8").Value = 1.142895050746318
$ws.Range("D8").Value = 0.1057680291380869
$ws.Range("E8").Value = 0.1765527611519611
$ws.Range("F8").Value = 2.666236566723924
$ws.Range("G8").Value = 1.921385949046609
$ws.Range("H8").Value = 1.644620933905742
$ws.Range("J8").Value = 0.2827583260077855
$ws.Range("K8").Value = 0.4835438149260938
$ws.Range("L8").Value = 0.180555874319694
$ws.Range("B9").Value = 1.196618512426767
$ws.Range("D9").Value = 0.105752147573245
$ws.Range("E9").Value = 0.1752229363972342
$ws.Range("F9").Value = 2.682938770788738
$ws.Range("G9").Value = 1.947003089062974
$ws.Range("H9").Value = 1.644420480881678
$ws.Range("J9").Value = 0.2800729976249485
$ws.Range("K9").Value = 0.5961864943121498
$ws.Range("L9").Value = 0.2000621912497991
$ws.Range("B10").Value = 1.238746992337212
$ws.Range("D10").Value = 0.1058505485410386
$ws.Range("E10").Value = 0.1743986352954572
$ws.Range("F10").Value = 2.700640065682109
$ws.Range("G10").Value = 1.970253513255756
$ws.Range("H10").Value = 1.6472650858164
$ws.Range("J10").Value = 0.2783214658242859
$ws.Range("K10").Value = 0.6799098820933978
$ws.Range("L10").Value = 0.2148573859723939
$ws.Range("B11").Value = 1.25848668746616
$ws.Range("D11").Value = 0.1059189407774639
$ws.Range("E11").Value = 0.1740565897039872
$ws.Range("F11").Value = 2.70987447861998
$ws.Range("G11").Value = 1.981797524639603
$ws.Range("H11").Value = 1.649209426202191
$ws.Range("J11").Value = 0.2775724865391016
$ws.Range("K11").Value = 0.7182078301418358
$ws.Range("L11").Value = 0.2216885365544954
$ws.Range("B12").Value = 1.266043967499826
$ws.Range("D12").Value = 0.1059482135259557
$ws.Range("E12").Value = 0.1739317848068005
$ws.Range("F12").Value = 2.713541397400249
$ws.Range("G12").Value = 1.986308347817925
$ws.Range("H12").Value = 1.650039227284367
$ws.Range("J12").Value = 0.2772957242516725
$ws.Range("K12").Value = 0.7327405624313315
$ws.Range("L12").Value = 0.2242897403403674
$ws.Range("B13").Value = 1.264412719191824
$ws.Range("D13").Value = 0.1059417594272496
$ws.Range("E13").Value = 0.1739584540802266
$ws.Range("F13").Value = 2.712744096721906
$ws.Range("G13").Value = 1.985330659218079
$ws.Range("H13").Value = 1.649856355253263
$ws.Range("J13").Value = 0.2773550250858037
$ws.Range("K13").Value = 0.7296093434804618
$ws.Range("L13").Value = 0.2237288857675139
$ws.Range("B14").Value = 1.259106782776939
$ws.Range("D14").Value = 0.1059212815765065
$ws.Range("E14").Value = 0.174046227415138
$ws.Range("F14").Value = 2.710172749801146
$ws.Range("G14").Value = 1.98216583835449
$ws.Range("H14").Value = 1.649275820081954
$ws.Range("J14").Value = 0.2775495797354601
$ws.Range("K14").Value = 0.7194028447260621
$ws.Range("L14").Value = 0.2219022511700075
$ws.Range("B15").Value = 1.2558674467885
$ws.Range("D15").Value = 0.1059091770087797
$ws.Range("E15").Value = 0.1741006053900023
$ws.Range("F15").Value = 2.708619872978147
$ws.Range("G15").Value = 1.980245450355028
$ws.Range("H15").Value = 1.64893240501948
$ws.Range("J15").Value = 0.2776696430293857
$ws.Range("K15").Value = 0.7131549840154605
$ws.Range("L15").Value = 0.2207852567452164
$ws.Range("B16").Value = 1.23746849503496
$ws.Range("D16").Value = 0.1058465523900516
$ws.Range("E16").Value = 0.1744216501595535
$ws.Range("F16").Value = 2.700060366532426
$ws.Range("G16").Value = 1.969518571556705
$ws.Range("H16").Value = 1.647151105184918
$ws.Range("J16").Value = 0.2783713740799403
$ws.Range("K16").Value = 0.6774112547521725
$ws.Range("L16").Value = 0.2144129735731894
$ws.Range("B17").Value = 1.226328357421011
$ws.Range("D17").Value = 0.1058141692170054
$ws.Range("E17").Value = 0.1746270247038151
$ws.Range("F17").Value = 2.69511218036493
$ws.Range("G17").Value = 1.96318590903013
$ws.Range("H17").Value = 1.646224901507139
$ws.Range("J17").Value = 0.2788140965281123
$ws.Range("K17").Value = 0.6555376415124385
$ws.Range("L17").Value = 0.2105295213633838
$ws.Range("B18").Value = 1.219975014702584
$ws.Range("D18").Value = 0.1057977687616543
$ws.Range("E18").Value = 0.1747482514463119
$ws.Range("F18").Value = 2.69237736584688
$ws.Range("G18").Value = 1.959634557939893
$ws.Range("H18").Value = 1.645753383904093
$ws.Range("J18").Value = 0.279073238881324
$ws.Range("K18").Value = 0.6429764682472126
$ws.Range("L18").Value = 0.2083053500804652
$ws.Range("B19").Value = 1.217833198980316
$ws.Range("D19").Value = 0.1057925988079838
$ws.Range("E19").Value = 0.1747898297631281
$ws.Range("F19").Value = 2.691470509862086
$ws.Range("G19").Value = 1.958447757767431
$ws.Range("H19").Value = 1.645604248875543
$ws.Range("J19").Value = 0.2791617533508113
$ws.Range("K19").Value = 0.6387269071833259
$ws.Range("L19").Value = 0.2075539162288891
$ws.Range("B20").Value = 1.227508640863277
$ws.Range("D20").Value = 0.105817386303773
$ws.Range("E20").Value = 0.1746048414305097
$ws.Range("F20").Value = 2.695627408581757
$ws.Range("G20").Value = 1.963850608703666
$ws.Range("H20").Value = 1.646317162383298
$ws.Range("J20").Value = 0.2787665023238741
$ws.Range("K20").Value = 0.6578640625168077
$ws.Range("L20").Value = 0.210941940203881
$ws.Range("B21").Value = 1.260663035134399
$ws.Range("D21").Value = 0.1059272050339644
$ws.Range("E21").Value = 0.1740203182724791
$ws.Range("F21").Value = 2.710923401084031
$ws.Range("G21").Value = 1.983091638760214
$ws.Range("H21").Value = 1.649443799124157
$ws.Range("J21").Value = 0.27749224825269
$ws.Range("K21").Value = 0.7223999256077605
$ws.Range("L21").Value = 0.2224383880656831
$ws.Range("B22").Value = 1.282810775302153
$ws.Range("D22").Value = 0.1060186339053608
$ws.Range("E22").Value = 0.1736658052998887
$ws.Range("F22").Value = 2.721911416469595
$ws.Range("G22").Value = 1.996479163484992
$ws.Range("H22").Value = 1.652032356235793
$ws.Range("J22").Value = 0.2766994280007093
$ws.Range("K22").Value = 0.7647533222240952
$ws.Range("L22").Value = 0.2300358468518482
$ws.Range("B23").Value = 1.270946383992055
$ws.Range("D23").Value = 0.1059680458441044
$ws.Range("E23").Value = 0.1738525037437508
$ws.Range("F23").Value = 2.715956185640977
$ws.Range("G23").Value = 1.98925956794028
$ws.Range("H23").Value = 1.650600913585095
$ws.Range("J23").Value = 0.2771189176176385
$ws.Range("K23").Value = 0.7421325652495057
$ws.Range("L23").Value = 0.2259732962224916
$ws.Range("B24").Value = 1.226974874982517
$ws.Range("D24").Value = 0.1058159249522959
$ws.Range("E24").Value = 0.1746148606627571
$ws.Range("F24").Value = 2.695394131368332
$ws.Range("G24").Value = 1.963549819576627
$ws.Range("H24").Value = 1.646275261314344
$ws.Range("J24").Value = 0.2787880052866551
$ws.Range("K24").Value = 0.6568122431467884
$ws.Range("L24").Value = 0.2107554592351022
$ws.Range("B25").Value = 1.181617071143648
$ws.Range("D25").Value = 0.1057370193641951
$ws.Range("E25").Value = 0.1755557973576609
$ws.Range("F25").Value = 2.677467521342052
$ws.Range("G25").Value = 1.939296674329711
$ws.Range("H25").Value = 1.643949454177232
$ws.Range("J25").Value = 0.2807604951045359
$ws.Range("K25").Value = 0.5655442008802538
$ws.Range("L25").Value = 0.1947036063300374
